$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.534.08"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "3.466.25"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.17"
$ws.Range("D5").NumberFormat = "general"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.81"
$ws.Range("D6").NumberFormat = "general"
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.467.46"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").NumberFormat = "general"
$ws.Range("E9").Value = "  +10.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("D10").NumberFormat = "general"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.127"
$ws.Range("D11").NumberFormat = "general"
$ws.Range("E11").Value = "  +5.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").NumberFormat = "general"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "4.056.96"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.10"
$ws.Range("D16").NumberFormat = "general"
$ws.Range("E16").Value = "  +7.76%  "
$ws.Range("D17").Value = "64.524.58"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "3.427.22"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").NumberFormat = "general"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.60"
$ws.Range("D20").NumberFormat = "general"
$ws.Range("E20").Value = "  +4.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.16"
$ws.Range("D21").NumberFormat = "general"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.26"
$ws.Range("D22").NumberFormat = "general"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("D23").NumberFormat = "general"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.33"
$ws.Range("D24").NumberFormat = "general"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("D26").NumberFormat = "general"
$ws.Range("E26").Value = "  +21.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("D27").NumberFormat = "general"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("D28").NumberFormat = "general"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.22"
$ws.Range("D30").NumberFormat = "general"
$ws.Range("E30").Value = "  +11.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").NumberFormat = "general"
$ws.Range("E31").Value = "  +10.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").NumberFormat = "general"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.79"
$ws.Range("D34").NumberFormat = "general"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.15"
$ws.Range("D36").NumberFormat = "general"
$ws.Range("E36").Value = "  +6.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("D37").NumberFormat = "general"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.09"
$ws.Range("D38").NumberFormat = "general"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").NumberFormat = "general"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0779"
$ws.Range("D40").NumberFormat = "general"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.66"
$ws.Range("D41").NumberFormat = "general"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "2.922.58"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.59"
$ws.Range("D43").NumberFormat = "general"
$ws.Range("E43").Value = "  +6.87%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.79"
$ws.Range("D45").NumberFormat = "general"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.776"
$ws.Range("D46").NumberFormat = "general"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.04"
$ws.Range("D47").NumberFormat = "general"
$ws.Range("E47").Value = "  +8.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("D48").NumberFormat = "general"
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("E49").Value = "  +17.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.876"
$ws.Range("D50").NumberFormat = "general"
$ws.Range("E50").Value = "  +7.92%  "
$ws.Range("E51").Value = "  +4.82%  "
